$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.688997030258179
$ws.Range("B1").Value = 3.287040233612061
$ws.Range("C1").Value = 1.600455045700073
$ws.Range("D1").Value = 1.324506759643555
$ws.Range("E1").Value = 1.398847579956055
